$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (single decimal
# point, e.g. "1.00", "6.37", "0.0000249"). Excel would normally parse these
# as numbers on assignment and drop the text formatting (e.g. "1.00" -> 1).
# Since the source data keeps these as text (prices formatted like
# "68.595.76" are text too), mark them as Text before writing the value so
# the literal string is preserved exactly.
$textRefs = @("D4","D5","D6","D11","D12","D13","D14","D19","D20","D21","D22","D23","D25","D28","D32","D33","D38","D42","D43","D44","D48","D49","D51")
$textRange = $ws.Range($textRefs[0])
foreach ($ref in $textRefs[1..($textRefs.Length - 1)]) {
    $textRange = $excel.Union($textRange, $ws.Range($ref))
}
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.595.76"
$ws.Range("D3").Value = "3.916.95"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "602.59"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").Value = "165.43"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").Value = "3.913.44"
$ws.Range("E7").Value = "  +3.05%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.49%  "
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "37.26"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "4.567.61"
$ws.Range("E15").Value = "  +3.05%  "
$ws.Range("D16").Value = "3.921.46"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "68.714.15"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "17.20"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "0.111"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").Value = "10.97"
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "485.48"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "0.724"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("E24").Value = "  +12.09%  "
$ws.Range("D25").Value = "84.60"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("D31").Value = "4.067.40"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("D32").Value = "2.38"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "7.75"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "3.864.65"
$ws.Range("E35").Value = "  +3.18%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Value = "5.92"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.315"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("D43").Value = "429.46"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").Value = "48.44"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "26.20"
$ws.Range("E48").Value = "  +7.26%  "
$ws.Range("D49").Value = "141.82"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "2.814.14"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "0.0352"
$ws.Range("E51").Value = "  +0.68%  "
